$wb = $excel.ActiveWorkbook

# Sheet1: ALC
$ws = $wb.Worksheets.Item(1)
$ws.Range("H17").Value = 950.6875
$ws.Range("J17").Value = 950.6875
$ws.Range("L17").Value = 2852.0625
$ws.Range("N17").Value = -3188.0625
$ws.Range("H40").Value = 9979.75
$ws.Range("I40").Value = 1570.7142
$ws.Range("K40").Value = 1570.7142
$ws.Range("M40").Value = -1395.7142
$ws.Range("H62").Value = 6254122.5
$ws.Range("I62").Value = 9618221
$ws.Range("K62").Value = 9618221
$ws.Range("M62").Value = -9617597
$ws.Range("H65").Value = 6254122.5
$ws.Range("I65").Value = 9618221
$ws.Range("K65").Value = 48091105
$ws.Range("M65").Value = -48087985
$ws.Range("H106").Value = 1828.4286
$ws.Range("I106").Value = 2084.5715
$ws.Range("K106").Value = 2084.5715
$ws.Range("M106").Value = -1453.5715
$ws.Range("H121").Value = 3781.182
$ws.Range("J121").Value = 3781.182
$ws.Range("L121").Value = 11343.546
$ws.Range("N121").Value = -14837.546
$ws.Range("H129").Value = 2633.973
$ws.Range("J129").Value = 2771.3438
$ws.Range("L129").Value = 8314.0314
$ws.Range("N129").Value = -18314.0314
$ws.Range("H137").Value = 12538.462
$ws.Range("I137").Value = 2247.6538
$ws.Range("J137").Value = 33120.08
$ws.Range("K137").Value = 6742.9614
$ws.Range("L137").Value = 99360.24000000001
$ws.Range("M137").Value = -4192.9614
$ws.Range("N137").Value = -104460.24
$ws.Range("H138").Value = 7090.26
$ws.Range("J138").Value = 11556.577
$ws.Range("L138").Value = 34669.731
$ws.Range("N138").Value = -44949.731

# Sheet2: ARM
$ws = $wb.Worksheets.Item(2)
$ws.Range("H32").Value = 2918.4238
$ws.Range("I32").Value = 2865.5
$ws.Range("K32").Value = 2865.5
$ws.Range("M32").Value = -2578.5
$ws.Range("H45").Value = 3519
$ws.Range("I45").Value = 3470.9
$ws.Range("K45").Value = 3470.9
$ws.Range("M45").Value = -3093.9
$ws.Range("H63").Value = 5715.231
$ws.Range("H66").Value = 5715.231
$ws.Range("H97").Value = 1710.931
$ws.Range("I97").Value = 1788.3077
$ws.Range("J97").Value = 1040.3334
$ws.Range("K97").Value = 1788.3077
$ws.Range("L97").Value = 1040.3334
$ws.Range("M97").Value = -1292.3077
$ws.Range("N97").Value = -2032.3334
$ws.Range("H122").Value = 38465144
$ws.Range("J122").Value = 6091.5
$ws.Range("L122").Value = 18274.5
$ws.Range("N122").Value = -23174.5
$ws.Range("H132").Value = 418820.47
$ws.Range("I132").Value = 569529.9399999999
$ws.Range("K132").Value = 1708589.82
$ws.Range("M132").Value = -1706059.82

# Sheet3: BSM
$ws = $wb.Worksheets.Item(3)
$ws.Range("H82").Value = 51049.777
$ws.Range("I82").Value = 11890.6
$ws.Range("J82").Value = 99998.75
$ws.Range("K82").Value = 11890.6
$ws.Range("L82").Value = 99998.75
$ws.Range("M82").Value = -11507.6
$ws.Range("N82").Value = -100764.75
$ws.Range("H85").Value = 51049.777
$ws.Range("I85").Value = 11890.6
$ws.Range("J85").Value = 99998.75
$ws.Range("K85").Value = 11890.6
$ws.Range("L85").Value = 99998.75
$ws.Range("M85").Value = -10564.6
$ws.Range("N85").Value = -102650.75
$ws.Range("H134").Value = 26354.69
$ws.Range("I134").Value = 1669.5428
$ws.Range("J134").Value = 112752.7
$ws.Range("K134").Value = 5008.6284
$ws.Range("L134").Value = 338258.1
$ws.Range("M134").Value = -2473.6284
$ws.Range("N134").Value = -343328.1

# Sheet4: CRP
$ws = $wb.Worksheets.Item(4)
$ws.Range("H99").Value = 5693.357
$ws.Range("I99").Value = 4384.857
$ws.Range("J99").Value = 7001.857
$ws.Range("K99").Value = 4384.857
$ws.Range("L99").Value = 7001.857
$ws.Range("M99").Value = -2886.857
$ws.Range("N99").Value = -9997.857
$ws.Range("H122").Value = 58743.668
$ws.Range("I122").Value = 333333
$ws.Range("K122").Value = 999999
$ws.Range("M122").Value = -997549
$ws.Range("H126").Value = 5693.357
$ws.Range("I126").Value = 4384.857
$ws.Range("J126").Value = 7001.857
$ws.Range("K126").Value = 13154.571
$ws.Range("L126").Value = 21005.571
$ws.Range("M126").Value = -10684.571
$ws.Range("N126").Value = -25945.571

# Sheet5: CUL
$ws = $wb.Worksheets.Item(5)
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").Value = $null
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").Value = $null
$ws.Range("H68").Value = 2138.5312
$ws.Range("J68").Value = 2159.8462
$ws.Range("L68").Value = 6479.5386
$ws.Range("N68").Value = -8101.5386
$ws.Range("H71").Value = 2138.5312
$ws.Range("J71").Value = 2159.8462
$ws.Range("L71").Value = 19438.6158
$ws.Range("N71").Value = -27550.6158
$ws.Range("H107").Value = 23084.898
$ws.Range("I107").Value = 618.8182
$ws.Range("J107").Value = 29588.236
$ws.Range("K107").Value = 1856.4546
$ws.Range("L107").Value = 88764.708
$ws.Range("M107").Value = 63.54539999999997
$ws.Range("N107").Value = -92604.708
$ws.Range("H131").Value = 46997.023
$ws.Range("I131").Value = 46363.316
$ws.Range("J131").Value = 47603.176
$ws.Range("K131").Value = 139089.948
$ws.Range("L131").Value = 142809.528
$ws.Range("M131").Value = -134049.948
$ws.Range("N131").Value = -152889.528

# Sheet6: GSM
$ws = $wb.Worksheets.Item(6)
$ws.Range("H97").Value = 1149.2069
$ws.Range("I97").Value = 1035.5264
$ws.Range("K97").Value = 1035.5264
$ws.Range("M97").Value = -539.5264

# Sheet7: LTW
$ws = $wb.Worksheets.Item(7)
$ws.Range("H3").Value = 10000
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").Value = $null
$ws.Range("H15").Value = 10000
$ws.Range("I15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("M15").Value = $null
$ws.Range("H46").Value = 3776.28
$ws.Range("J46").Value = 5240.4
$ws.Range("L46").Value = 5240.4
$ws.Range("N46").Value = -5616.4
$ws.Range("H122").Value = 374146.94
$ws.Range("J122").Value = 855461.7
$ws.Range("L122").Value = 2566385.1
$ws.Range("N122").Value = -2571285.1
$ws.Range("H132").Value = 4073.2964
$ws.Range("I132").Value = 3583
$ws.Range("J132").Value = 5237.75
$ws.Range("K132").Value = 10749
$ws.Range("L132").Value = 15713.25
$ws.Range("M132").Value = -8219
$ws.Range("N132").Value = -20773.25
$ws.Range("H135").Value = 74999.5
$ws.Range("J135").Value = 74999.5
$ws.Range("L135").Value = 74999.5
$ws.Range("N135").Value = -85139.5

# Sheet8: WVR
$ws = $wb.Worksheets.Item(8)
$ws.Range("H19").Value = 2026.25
$ws.Range("J19").Value = 2201.6667
$ws.Range("L19").Value = 2201.6667
$ws.Range("N19").Value = -2549.6667
$ws.Range("H132").Value = 49910.168
$ws.Range("I132").Value = 3531.7646
$ws.Range("M132").Value = -8065.293799999999
